$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "67.721.74"
$ws.Range("E2").Value = "  +0.47%  "

# Row 3
$ws.Range("D3").Value = "2.492.31"
$ws.Range("E3").Value = "  -2.46%  "

# Row 4
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
$ws.Range("D5").Value = "'591.55"
$ws.Range("E5").Value = "  -0.31%  "

# Row 6
$ws.Range("D6").Value = "'173.16"
$ws.Range("E6").Value = "  -0.29%  "

# Row 7
$ws.Range("E7").Value = "  -0.02%  "

# Row 8
$ws.Range("E8").Value = "  -0.96%  "

# Row 9
$ws.Range("D9").Value = "2.491.37"
$ws.Range("E9").Value = "  -2.45%  "

# Row 10
$ws.Range("D10").Value = "'0.140"
$ws.Range("E10").Value = "  -0.27%  "

# Row 11
$ws.Range("E11").Value = "  +1.83%  "

# Row 12
$ws.Range("E12").Value = "  -1.91%  "

# Row 13
$ws.Range("E13").Value = "  -3.16%  "

# Row 14
$ws.Range("D14").Value = "'26.23"
$ws.Range("E14").Value = "  -3.38%  "

# Row 15
$ws.Range("D15").Value = "2.944.50"
$ws.Range("E15").Value = "  -2.01%  "

# Row 16
$ws.Range("D16").Value = "'0.0000176"
$ws.Range("E16").Value = "  -1.52%  "

# Row 17
$ws.Range("D17").Value = "67.605.78"
$ws.Range("E17").Value = "  +0.42%  "

# Row 18
$ws.Range("D18").Value = "2.470.26"
$ws.Range("E18").Value = "  -3.04%  "

# Row 19
$ws.Range("E19").Value = "  +2.85%  "

# Row 20
$ws.Range("D20").Value = "'7.97"
$ws.Range("E20").Value = "  +0.26%  "

# Row 21
$ws.Range("D21").Value = "'363.90"
$ws.Range("E21").Value = "  +1.96%  "

# Row 22
$ws.Range("E22").Value = "  -2.67%  "

# Row 23
$ws.Range("D23").Value = "'4.55"
$ws.Range("E23").Value = "  -3.19%  "

# Row 24
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Value = "'71.19"
$ws.Range("E24").Value = "  +1.37%  "

# Row 25
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").Value = "'0.998"
$ws.Range("E25").Value = "  -0.12%  "

# Row 26
$ws.Range("E26").Value = "  -5.92%  "

# Row 27
$ws.Range("D27").Value = "'9.80"
$ws.Range("E27").Value = "  -2.51%  "

# Row 28
$ws.Range("D28").Value = "'0.998"
$ws.Range("E28").Value = "  -0.12%  "

# Row 29
$ws.Range("D29").Value = "2.610.38"
$ws.Range("E29").Value = "  -2.94%  "

# Row 30
$ws.Range("D30").Value = "0.0₃0962"
$ws.Range("E30").Value = "  -3.88%  "

# Row 31
$ws.Range("D31").Value = "'531.08"
$ws.Range("E31").Value = "  -1.10%  "

# Row 32
$ws.Range("D32").Value = "'8.23"
$ws.Range("E32").Value = "  -0.58%  "

# Row 33
$ws.Range("D33").Value = "'1.86"
$ws.Range("E33").Value = "  +0.00%  "

# Row 34
$ws.Range("E34").Value = "  -4.66%  "

# Row 35
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  +0.01%  "

# Row 36
$ws.Range("E36").Value = "  -4.43%  "

# Row 37
$ws.Range("D37").Value = "'159.11"
$ws.Range("E37").Value = "  +1.13%  "

# Row 38
$ws.Range("D38").Value = "'1.42"
$ws.Range("E38").Value = "  -4.07%  "

# Row 39
$ws.Range("E39").Value = "  -1.25%  "

# Row 40
$ws.Range("D40").Value = "'18.64"
$ws.Range("E40").Value = "  +0.98%  "

# Row 41
$ws.Range("E41").Value = "  -1.97%  "

# Row 42
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D42").Value = "'5.11"
$ws.Range("E42").Value = "  -2.28%  "

# Row 43
$ws.Range("B43").Value = "PolygonEcosystemToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D43").Value = "'0.348"
$ws.Range("E43").Value = "  -2.97%  "

# Row 44
$ws.Range("D44").Value = "'0.999"
$ws.Range("E44").Value = "  -0.15%  "

# Row 45
$ws.Range("E45").Value = "  -3.31%  "

# Row 46
$ws.Range("D46").Value = "'144.71"
$ws.Range("E46").Value = "  -4.44%  "

# Row 47
$ws.Range("D47").Value = "'3.68"
$ws.Range("E47").Value = "  -1.51%  "

# Row 48
$ws.Range("D48").Value = "'0.546"
$ws.Range("E48").Value = "  -3.66%  "

# Row 49
$ws.Range("D49").Value = "0.0₆0271"
$ws.Range("E49").Value = "  -4.95%  "

# Row 50
$ws.Range("D50").Value = "'1.69"
$ws.Range("E50").Value = "  -2.21%  "

# Row 51
$ws.Range("E51").Value = "  -1.90%  "

